# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Mon Nov 25 17:50:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.076.82"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.437.54"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'238.80"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").Value = "'641.72"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +7.55%  "
$ws.Range("D8").Value = "'0.401"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").Value = "3.432.29"
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("D13").Value = "'41.40"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").Value = "'6.08"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "94.892.77"
$ws.Range("D16").Value = "4.080.75"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "'8.38"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "3.453.98"
$ws.Range("E19").Value = "  +4.32%  "
$ws.Range("D20").Value = "'17.88"
$ws.Range("E20").Value = "  +6.91%  "
$ws.Range("D21").Value = "'11.52"
$ws.Range("E21").Value = "  +10.88%  "
$ws.Range("D22").Value = "'0.511"
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").Value = "'499.86"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "'6.55"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'91.34"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.618.84"
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'11.98"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  +8.64%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +11.91%  "
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'0.182"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("E35").Value = "  +10.96%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.564"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("D39").Value = "'1.43"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "'511.07"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("D42").Value = "'0.149"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'0.907"
$ws.Range("E43").Value = "  +11.07%  "
$ws.Range("D44").Value = "'24.12"
$ws.Range("E45").Value = "  +5.28%  "
$ws.Range("D46").Value = "'0.0414"
$ws.Range("D47").Value = "'5.50"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'53.43"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.13"
$ws.Range("E50").Value = "  +9.72%  "
$ws.Range("E51").Value = "  +2.53%  "
